# Data drive test case 1
# Adds a new worksheet "ValidLogin" (placed right after "TC1") with
# username/password test data, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet right after the existing TC1 sheet and rename it.
$tc1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $tc1)
$newSheet.Name = "ValidLogin"

# Populate the new sheet with the data-driven login values.
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "pointofsale"

# Make the new "ValidLogin" sheet the active tab.
$newSheet.Activate()

# Select B3 as the active cell, matching the saved selection state.
$newSheet.Range("B3").Select()

# Match the saved zoom level for the new sheet.
$excel.ActiveWindow.Zoom = 160
